$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 52, shifting existing rows 52:55 down to 53:56
$ws.Rows.Item(52).EntireRow.Insert()

# Populate the newly inserted row 52 with the new weekly data point
$ws.Cells.Item(52, 1).Value = 1
$ws.Cells.Item(52, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(52, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(52, 4).Value = 44617
$ws.Cells.Item(52, 5).Value = 15
$ws.Cells.Item(52, 6).Value = 100112009
$ws.Cells.Item(52, 7).Value = "Acelga"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Segunda"
$ws.Cells.Item(52, 10).Value = 250
$ws.Cells.Item(52, 11).Value = 1000
$ws.Cells.Item(52, 12).Value = 1200
$ws.Cells.Item(52, 13).Value = 1100
$ws.Cells.Item(52, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(52, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(52, 16).Value = 367
$ws.Cells.Item(52, 17).Value = 3
$ws.Cells.Item(52, 18).Value = "Hortaliza"

# Keep D52's number format consistent with the other date cells in column D
$ws.Cells.Item(52, 4).NumberFormat = $ws.Cells.Item(51, 4).NumberFormat
